# Add a new "Week 17" row (spreadsheet row 18) to the Comprehension scores sheet,
# mirroring the pattern used by the existing rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A18: Week number
$ws.Range("A18").Value2 = 17

# B18: Total time for the week (same style/number format as B17, i.e. [h]:mm:ss)
$ws.Range("B18").Value2 = 1.6707291666666666
$ws.Range("B18").NumberFormat = $ws.Range("B17").NumberFormat

# C18: Running-total formula identical in shape to the one used in C17, extended one row
$ws.Range("C18").Formula = "=SUM(B2:B18)+1.2708333333"
$ws.Range("C18").NumberFormat = $ws.Range("C17").NumberFormat

# D18: Comprehension scores text for the new week
$ws.Range("D18").Value2 = "Élite (Subtitled, Spanish, New):33; Harry Potter book 5 (Text-only, English, Familiar):37; [The Last Airbender: Una HORRIBLE Película de Avatar | Tortura Audiovisual | LA ZONA CERO](https://youtu.be/_k91HM04XNg) (Audiovisual, Spanish, New):38; La venganza de las Juanas  (Audiovisual, Spanish, New):32; Velvet (Audiovisual, Spanish, New):38; Siempre bruja (Subtitled, Spanish, New):41; [El teorema de Sheldon Cooper: El número 73 es único](https://youtu.be/R7hTUxzbH48) (Audiovisual, Spanish, New):33; Historia de un crimen: Colmenares (Subtitled, Spanish, New):39;"

# Move the active selection to C19, matching where Excel would land after entering this row
$ws.Range("C19").Select()
